# Telnet and SSH tested on 3 devices
# Mark the SSH access (C2) result as OK for the Cisco IOS device (column C),
# matching the style/value already used in the neighboring OK cells.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the "OK" formatting (green, centered) from B2 onto C2, then set its value.
$ws.Range("B2").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C2").Value = "OK"

# Mark the Telnet access (C4) result as tested ("X") for the Cisco IOS device.
$ws.Range("C4").Value = "X"

# Rows 2-4 re-wrap to the compact (13.8pt) row height used elsewhere in the sheet.
$ws.Rows.Item(2).RowHeight = 13.8
$ws.Rows.Item(3).RowHeight = 13.8
$ws.Rows.Item(4).RowHeight = 13.8

# Reflect the cells the author was reviewing/editing when saving the file.
$selectionResult = $ws.Range("C2:C4").Select()
